# Update "想去人数" (want-to-go count) figures in column F of the
# "展览" (Exhibitions) and "全部类型" (All types) sheets to reflect the
# refreshed scrape output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$ws1.Range("F2").Value  = 1253    # 昆山·创世次元动漫游戏嘉年华
$ws1.Range("F3").Value  = 17116   # 苏州·ICAN summer World动漫品牌夏游节
$ws1.Range("F8").Value  = 1049    # 苏州·幻想物语次元嘉年华（免费展）
$ws1.Range("F9").Value  = 404     # 苏州·排球少年only-茶歇
$ws1.Range("F12").Value = 11904   # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$ws1.Range("F14").Value = 57      # 【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场
$ws1.Range("F15").Value = 11604   # 苏州·I COME ACG动漫品牌博览会
$ws1.Range("F16").Value = 4719    # 苏州·理想乡动漫游戏展-两馆全开+三馆间通道
$ws1.Range("F17").Value = 509     # 苏州·第四届-OCG国朝动漫游戏嘉年华
$ws1.Range("F19").Value = 416     # 苏州·明日方舟ONLY#2024~佑桑柔
$ws1.Range("F25").Value = 5215    # 苏州·星部落动漫嘉年华

# --- Sheet "全部类型" ---
$ws4.Range("F2").Value  = 1253    # 昆山·创世次元动漫游戏嘉年华
$ws4.Range("F3").Value  = 17116   # 苏州·ICAN summer World动漫品牌夏游节
$ws4.Range("F8").Value  = 1049    # 苏州·幻想物语次元嘉年华（免费展）
$ws4.Range("F9").Value  = 404     # 苏州·排球少年only-茶歇
$ws4.Range("F14").Value = 11904   # 苏州·Good jump ACG中秋嘉年华动漫国潮文化节
$ws4.Range("F16").Value = 57      # 【大会员提前抢】苏州·OCG国潮动漫游戏嘉年华CV杨天翔内场
$ws4.Range("F17").Value = 11605   # 苏州·I COME ACG动漫品牌博览会
$ws4.Range("F18").Value = 4719    # 苏州·理想乡动漫游戏展-两馆全开+三馆间通道
$ws4.Range("F19").Value = 509     # 苏州·第四届-OCG国朝动漫游戏嘉年华
$ws4.Range("F21").Value = 416     # 苏州·明日方舟ONLY#2024~佑桑柔
$ws4.Range("F27").Value = 5215    # 苏州·星部落动漫嘉年华

$wb.Save()
